# FL21_DATA317_Electrical_Usage.xlsx - "Add files via upload" revision
#
# The sheet goes from a 5-column (day, Grant, Hoyum, Memorial, Library),
# 10-data-row table (plus some orphaned leftover values in A12:A21) to a
# 6-column table that prepends a calendar "Date" column and carries the
# data out to 18 days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Date" column in front of the old data -----------------
# This shifts the existing day/Grant/Hoyum/Memorial/Library columns from
# A:E to B:F (carrying their values + styles with them) and leaves a blank
# column A to populate.
$ws.Columns("A:A").Insert()

# The new column A should be the same width as the (old) day column used to
# be before the insert.
$ws.Columns("A:A").ColumnWidth = 13.17

# --- Header row --------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Date"
$ws.Cells.Item(1, 2).Value = "day"
$ws.Cells.Item(1, 3).Value = "Grant"
$ws.Cells.Item(1, 4).Value = "Hoyum"
$ws.Cells.Item(1, 5).Value = "Memorial"
$ws.Cells.Item(1, 6).Value = "Library"

# --- Data rows (18 days of readings) ------------------------------------
$ws.Cells.Item(2, 1).Value = 44445
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 4282.5
$ws.Cells.Item(2, 4).Value = 469.06
$ws.Cells.Item(2, 5).Value = 1017.5
$ws.Cells.Item(2, 6).Value = 1548.72
$ws.Cells.Item(3, 1).Value = 44446
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 4417.91
$ws.Cells.Item(3, 4).Value = 480.63
$ws.Cells.Item(3, 5).Value = 1041.88
$ws.Cells.Item(3, 6).Value = 1546.31
$ws.Cells.Item(4, 1).Value = 44447
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 3962.91
$ws.Cells.Item(4, 4).Value = 437.31
$ws.Cells.Item(4, 5).Value = 980.44
$ws.Cells.Item(4, 6).Value = 1497.8
$ws.Cells.Item(5, 1).Value = 44448
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 4086.81
$ws.Cells.Item(5, 4).Value = 413.44
$ws.Cells.Item(5, 5).Value = 983.56
$ws.Cells.Item(5, 6).Value = 1533.56
$ws.Cells.Item(6, 1).Value = 44449
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 4187.53
$ws.Cells.Item(6, 4).Value = 450.11
$ws.Cells.Item(6, 5).Value = 1005.84
$ws.Cells.Item(6, 6).Value = 1531.59
$ws.Cells.Item(7, 1).Value = 44450
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 3622.53
$ws.Cells.Item(7, 4).Value = 411.06
$ws.Cells.Item(7, 5).Value = 803.94
$ws.Cells.Item(7, 6).Value = 1220.14
$ws.Cells.Item(8, 1).Value = 44451
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(8, 3).Value = 3477.22
$ws.Cells.Item(8, 4).Value = 409.88
$ws.Cells.Item(8, 5).Value = 701.63
$ws.Cells.Item(8, 6).Value = 1301.36
$ws.Cells.Item(9, 1).Value = 44452
$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(9, 3).Value = 2892.28
$ws.Cells.Item(9, 4).Value = 422.19
$ws.Cells.Item(9, 5).Value = 949.63
$ws.Cells.Item(9, 6).Value = 1510.84
$ws.Cells.Item(10, 1).Value = 44453
$ws.Cells.Item(10, 2).Value = 9
$ws.Cells.Item(10, 3).Value = 3504.22
$ws.Cells.Item(10, 4).Value = 415
$ws.Cells.Item(10, 5).Value = 1013.75
$ws.Cells.Item(10, 6).Value = 1533.88
$ws.Cells.Item(11, 1).Value = 44454
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(11, 3).Value = 3454.34
$ws.Cells.Item(11, 4).Value = 401.5
$ws.Cells.Item(11, 5).Value = 1024.75
$ws.Cells.Item(11, 6).Value = 1544.22
$ws.Cells.Item(12, 1).Value = 44455
$ws.Cells.Item(12, 2).Value = 11
$ws.Cells.Item(12, 3).Value = 3728.78
$ws.Cells.Item(12, 4).Value = 411.56
$ws.Cells.Item(12, 5).Value = 967.5
$ws.Cells.Item(12, 6).Value = 1583.58
$ws.Cells.Item(13, 1).Value = 44456
$ws.Cells.Item(13, 2).Value = 12
$ws.Cells.Item(13, 3).Value = 3026.16
$ws.Cells.Item(13, 4).Value = 404.75
$ws.Cells.Item(13, 5).Value = 980.13
$ws.Cells.Item(13, 6).Value = 1326.05
$ws.Cells.Item(14, 1).Value = 44457
$ws.Cells.Item(14, 2).Value = 13
$ws.Cells.Item(14, 3).Value = 3501.97
$ws.Cells.Item(14, 4).Value = 367.81
$ws.Cells.Item(14, 5).Value = 739.25
$ws.Cells.Item(14, 6).Value = 1183.16
$ws.Cells.Item(15, 1).Value = 44458
$ws.Cells.Item(15, 2).Value = 14
$ws.Cells.Item(15, 3).Value = 4272.03
$ws.Cells.Item(15, 4).Value = 427.44
$ws.Cells.Item(15, 5).Value = 704.88
$ws.Cells.Item(15, 6).Value = 1379.98
$ws.Cells.Item(16, 1).Value = 44459
$ws.Cells.Item(16, 2).Value = 15
$ws.Cells.Item(16, 3).Value = 2494.72
$ws.Cells.Item(16, 4).Value = 439.44
$ws.Cells.Item(16, 5).Value = 990.5
$ws.Cells.Item(16, 6).Value = 1510.83
$ws.Cells.Item(17, 1).Value = 44460
$ws.Cells.Item(17, 2).Value = 16
$ws.Cells.Item(17, 3).Value = 2543.16
$ws.Cells.Item(17, 4).Value = 407.88
$ws.Cells.Item(17, 5).Value = 989.13
$ws.Cells.Item(17, 6).Value = 1513.8
$ws.Cells.Item(18, 1).Value = 44461
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(18, 3).Value = 2960.94
$ws.Cells.Item(18, 4).Value = 404.69
$ws.Cells.Item(18, 5).Value = 944.75
$ws.Cells.Item(18, 6).Value = 1508.41
$ws.Cells.Item(19, 1).Value = 44462
$ws.Cells.Item(19, 2).Value = 18
$ws.Cells.Item(19, 3).Value = 3227.66
$ws.Cells.Item(19, 4).Value = 405.69
$ws.Cells.Item(19, 5).Value = 813.19
$ws.Cells.Item(19, 6).Value = 1554.05

# --- Remove the old leftover/orphaned rows (20 & 21) ------------------
$ws.Rows("20:21").Delete()

# --- Formatting ----------------------------------------------------------
# Give the new Date column a short date display ("6-Sep" style).
$ws.Range("A2:A19").NumberFormat = "d-mmm"

# --- Selection, matching the saved workbook's last active cell -----------
$ws.Range("C21").Select()
